$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'22.383.24"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = "'1.568.03"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'1.003"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').Value = "'291.54"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.72%  '
$ws.Range('D7').Value = "'0.3763"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.28%  '
$ws.Range('D8').Value = "'49.65"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.61%  '
$ws.Range('E9').Value = '  +0.86%  '
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('D11').Value = "'1.141"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.78%  '
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').Value = "'21.07"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.86%  '
$ws.Range('D14').Value = "'5.986"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').Value = "'6.958"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.78%  '
$ws.Range('D16').Value = "'1.565.35"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('D18').Value = "'90.02"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('D19').Value = "'0.06737"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D21').Value = "'16.59"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('E22').Value = '  -0.63%  '
$ws.Range('D23').Value = "'11.94"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.26%  '
$ws.Range('D24').Value = "'22.376.04"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.08%  '
$ws.Range('D25').Value = "'2.385"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.48%  '
$ws.Range('D26').Value = "'2.688"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.32%  '
$ws.Range('E27').Value = '  +0.27%  '
$ws.Range('D28').Value = "'147.49"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.88%  '
$ws.Range('D29').Value = "'5.033"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.35%  '
$ws.Range('D30').Value = "'126.37"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.84%  '
$ws.Range('D31').Value = "'1.748.44"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.24%  '
$ws.Range('D32').Value = "'2.017"
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').Value = "'0.9973"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.79%  '
$ws.Range('D34').Value = "'6.077"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.74%  '
$ws.Range('D35').Value = "'10.12"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.77%  '
$ws.Range('D36').Value = "'0.08451"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.22%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = "'0.02514"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.84%  '
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = "'1.379"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.72%  '
$ws.Range('E39').Value = '  -0.98%  '
$ws.Range('D40').Value = "'0.06510"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.26%  '
$ws.Range('D41').Value = "'5.393"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.53%  '
$ws.Range('D42').Value = "'11.36"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.85%  '
$ws.Range('D43').Value = "'0.6325"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.40%  '
$ws.Range('D44').Value = "'1.003"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('D45').Value = "'14.01"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.43%  '
$ws.Range('D46').Value = "'3.802"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.50%  '
$ws.Range('D47').Value = "'0.5930"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.85%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = "'2.076"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.41%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D49').Value = "'1.276"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.34%  '
$ws.Range('D50').Value = "'124.12"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('D51').Value = "'0.07312"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.38%  '
